$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=-105.5674577863067;  C=0.09916957484935195; D=297.092816532},
    @{Row=3;  B=-105.84488407880603; C=0.09920532363664455; D=229.387371966},
    @{Row=4;  B=-105.4938270439315;  C=0.09013521694606508; D=238.119057973},
    @{Row=5;  B=-104.60548337657738; C=0.09852701461338081; D=304.863617716},
    @{Row=6;  B=-103.49903687190475; C=0.08675010732014433; D=219.10969867},
    @{Row=7;  B=-104.63902673475309; C=0.08112826102326275; D=219.021526803},
    @{Row=8;  B=-102.87894895186048; C=0.09775684673713221; D=213.943059503},
    @{Row=9;  B=-105.10124785893194; C=0.09898214076070072; D=201.424226699},
    @{Row=10; B=-105.9960833017315;  C=0.07092340862122525; D=251.647956802},
    @{Row=11; B=-102.78699673530178; C=0.09590964376430258; D=228.107253047}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}
